$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, pushing the old "Test_Cell" row down to row 3.
$ws.Rows.Item(2).Insert()

# Fill in the new data set (order chosen to reproduce the shared-string table
# layout of the target workbook).
$ws.Cells.Item(3,1).Value = "Cell2"
$ws.Cells.Item(4,1).Value = "Cell3"
$ws.Cells.Item(1,2).Value = "Chemistry"
$ws.Cells.Item(2,2).Value = "NMC622"
$ws.Cells.Item(4,2).Value = "LFP"
$ws.Cells.Item(1,3).Value = "Nominal Capacity [Ah]"
$ws.Cells.Item(2,1).Value = "Cell1"
$ws.Cells.Item(1,4).Value = "Start date"
$ws.Cells.Item(3,2).Value = "NMC811"

$ws.Cells.Item(2,3).Value = 5
$ws.Cells.Item(3,3).Value = 3
$ws.Cells.Item(4,3).Value = 2.5

$ws.Cells.Item(2,4).Value = 45371.377349537041
$ws.Cells.Item(2,4).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(3,4).Value = 45371.376655092594
$ws.Cells.Item(3,4).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(4,4).Value = 45371.377349537041
$ws.Cells.Item(4,4).NumberFormat = "m/d/yy h:mm"

# Column D best-fit width (closest achievable value via ColumnWidth, which is
# quantized to pixel widths by the host).
$ws.Columns.Item(4).ColumnWidth = 14.498697916666666

# Page setup: Letter-ish A4 (paperSize 9), portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Final selection, matching the saved view state.
$null = $ws.Range("E4").Select()

Write-Output "done"
